$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column F (dSF) values per repulled data
$ws.Range("F2").Value = -2
$ws.Range("F3").Value = -5
$ws.Range("F5").Value = 0
$ws.Range("F13").Value = -2
$ws.Range("F14").Value = -3
$ws.Range("F16").Value = 0
